$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "08302023"
$ws.Range("F3").Value = "08302023"
$ws.Range("E6").Select()
